$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Datatype aaa"
$ws.Range("B3").Value = "String"
$ws.Range("C3").Value = "aaaaaaa"

$ws.Range("B2:C2").Merge()
$ws.Range("B2:C3").Borders.LineStyle = 1
